$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.427.36"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "3.285.44"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'600.60"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").Value = "'138.77"
$ws.Range("E6").Value = "  -2.51%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").Value = "3.285.49"
$ws.Range("D9").Value = "'0.514"
$ws.Range("E9").Value = "  -1.40%  "
$ws.Range("D10").Value = "'0.148"
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("D11").Value = "'5.45"
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("D12").Value = "'0.463"
$ws.Range("E12").Value = "  -1.43%  "
$ws.Range("D13").Value = "'0.0000243"
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").Value = "'34.28"
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "3.814.46"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("E16").Value = "  +0.56%  "
$ws.Range("D17").Value = "3.274.46"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "63.400.60"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").Value = "'6.81"
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("D20").Value = "'473.87"
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("D21").Value = "'13.96"
$ws.Range("E21").Value = "  -2.25%  "
$ws.Range("D22").Value = "'0.731"
$ws.Range("E22").Value = "  -1.01%  "
$ws.Range("D23").Value = "'7.88"
$ws.Range("E23").Value = "  -2.16%  "
$ws.Range("D24").Value = "'13.77"
$ws.Range("E24").Value = "  +1.70%  "
$ws.Range("D25").Value = "'85.22"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  -0.82%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'8.06"
$ws.Range("E29").Value = "  -1.60%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'7.03"
$ws.Range("E30").Value = "  -4.00%  "
$ws.Range("D31").Value = "'2.12"
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("D32").Value = "'28.49"
$ws.Range("E32").Value = "  +2.34%  "
$ws.Range("D33").Value = "'0.104"
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("D34").Value = "'2.50"
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("D36").Value = "'5.98"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").Value = "'52.03"
$ws.Range("D38").Value = "0.0₃0726"
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("D39").Value = "'0.0401"
$ws.Range("E39").Value = "  +1.41%  "
$ws.Range("D40").Value = "3.087.89"
$ws.Range("E40").Value = "  +3.25%  "
$ws.Range("D41").Value = "'426.83"
$ws.Range("E41").Value = "  +1.34%  "
$ws.Range("D42").Value = "'0.118"
$ws.Range("E42").Value = "  +6.18%  "
$ws.Range("D43").Value = "'8.24"
$ws.Range("D44").Value = "'2.70"
$ws.Range("E44").Value = "  -3.34%  "
$ws.Range("D45").Value = "'0.259"
$ws.Range("E45").Value = "  -2.98%  "
$ws.Range("D46").Value = "'2.19"
$ws.Range("E46").Value = "  -1.75%  "
$ws.Range("D47").Value = "'36.76"
$ws.Range("E47").Value = "  +8.61%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'26.12"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").Value = "'0.998"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'126.29"
$ws.Range("E50").Value = "  +3.03%  "
$ws.Range("E51").Value = "  -1.21%  "
